$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New soil-sample rows (SS-6 .. SS-10), continuing the existing table
# that currently ends at row 8 (SS-5).
$newRows = @(
    @("SS-6",  6,    6.5,  11.3, 1.8, 99, 20, 38),
    @("SS-7",  6.25, 6.5,  14.2, 1.8, 99, 22, 40),
    @("SS-8",  6.5,  6.73, 9.8,  1.8, 98, 23, 100),
    @("SS-9",  7.5,  7.63, 6.6,  1.8, 84, 15, 100),
    @("SS-10", 9,    9.1,  10,   1.8, 84, 15, 100)
)

$r = 9
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]   # A - SampleNo
    $ws.Cells.Item($r, 2).Value = $row[1]   # B - From
    $ws.Cells.Item($r, 3).Value = $row[2]   # C - To
    $ws.Cells.Item($r, 4).Value = $row[3]   # D - Wn %
    $ws.Cells.Item($r, 5).Value = $row[4]   # E - y (t/cu.m.)
    $ws.Cells.Item($r, 6).Value = $row[5]   # F - #4 (%)
    $ws.Cells.Item($r, 7).Value = $row[6]   # G - #200 (%)
    $ws.Cells.Item($r, 14).Value = $row[7]  # N - N
    $r++
}

# Match the author's final cursor position/selection after entering the data.
$ws.Range("N14").Select()

# Restore the window chrome size/position recorded in the saved file
# (best-effort; harmless if the host does not persist window geometry).
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Top = -120
$win.Left = -120
$win.Width = 29040
$win.Height = 15840
